# RubricaAndamento.xlsx — update "Escopo" progress rubric
# Commit: "Implementação dos dtos Entidades e ajustes para gravar as ordem de serviço"
#
# The F column holds completion percentages (0..1) for each rubric sub-item;
# row 26 (F26) keeps an AVERAGE(F2:F25) formula that recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escopo")
$ws.Activate()

# --- Progress updates -------------------------------------------------
$ws.Range("F4").Value  = 1        # "O aluno criou as migrações do banco de dados?" 0.8 -> 1
$ws.Range("F6").Value  = 0.6      # "O aluno criou as entidades, value objects?"    0.35 -> 0.6
$ws.Range("F7").Value  = 0.7      # "O aluno utilizou validações de entidades?"     0.5 -> 0.7
$ws.Range("F9").Value  = 0.66     # "O aluno utilizou CQRS..."                      0.15 -> 0.66
$ws.Range("F10").Value = 0.4      # "O aluno criou API REST leitura?"               0.15 -> 0.4
$ws.Range("F11").Value = 0.4      # "O aluno criou API REST escrita?"               0.15 -> 0.4
$ws.Range("F12").Value = 1        # "O aluno utilizou DTO..."                       0.25 -> 1
$ws.Range("F17").Value = 0.5      # "O aluno utilizou Tables..."                    0 -> 0.5
$ws.Range("F19").Value = 0.5      # "O aluno mapeou corretamente..."                0 -> 0.5

# --- Selection / view state -------------------------------------------
$ws.Range("F13").Select()

$wb.Save()
